$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format column D as Text so numeric-looking price strings
# (e.g. "4.83", "1.00") are preserved verbatim instead of being
# auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "68.524.72"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "2.460.45"
$ws.Range("E3").Value = "  -1.49%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "561.46"
$ws.Range("E5").Value = "  -2.26%  "
$ws.Range("D6").Value = "164.36"
$ws.Range("E6").Value = "  -1.16%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "0.504"
$ws.Range("E8").Value = "  -1.42%  "
$ws.Range("D9").Value = "2.459.43"
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("E10").Value = "  -4.94%  "
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").Value = "0.340"
$ws.Range("E12").Value = "  -4.67%  "
$ws.Range("D13").Value = "4.83"
$ws.Range("E13").Value = "  -1.98%  "
$ws.Range("D14").Value = "2.904.81"
$ws.Range("E14").Value = "  -1.73%  "
$ws.Range("D15").Value = "68.445.66"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "0.0000171"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").Value = "23.49"
$ws.Range("E17").Value = "  -4.75%  "
$ws.Range("D18").Value = "2.479.62"
$ws.Range("E18").Value = "  -0.78%  "
$ws.Range("D19").Value = "11.03"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("D20").Value = "7.22"
$ws.Range("E20").Value = "  -2.76%  "
$ws.Range("D21").Value = "344.43"
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("D22").Value = "3.80"
$ws.Range("E22").Value = "  -2.57%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "1.88"
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").Value = "67.85"
$ws.Range("E25").Value = "  -3.98%  "
$ws.Range("D26").Value = "1.12"
$ws.Range("E26").Value = "  +11.94%  "
$ws.Range("D27").Value = "3.74"
$ws.Range("E27").Value = "  -4.73%  "
$ws.Range("D28").Value = "2.583.43"
$ws.Range("E28").Value = "  -1.72%  "
$ws.Range("D29").Value = "8.21"
$ws.Range("E29").Value = "  -6.06%  "
$ws.Range("D30").Value = "0.0₃0840"
$ws.Range("E30").Value = "  -5.43%  "
$ws.Range("D31").Value = "7.29"
$ws.Range("E31").Value = "  -6.93%  "
$ws.Range("D32").Value = "3.45"
$ws.Range("E32").Value = "  +134.62%  "
$ws.Range("D33").Value = "435.02"
$ws.Range("E33").Value = "  -4.84%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.18"
$ws.Range("E34").Value = "  -2.55%  "
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "1.68"
$ws.Range("E36").Value = "  -2.69%  "
$ws.Range("D37").Value = "157.92"
$ws.Range("E37").Value = "  +1.14%  "
$ws.Range("D38").Value = "19.00"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.110"
$ws.Range("E39").Value = "  -4.17%  "
$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.01%  "
$ws.Range("D41").Value = "17.90"
$ws.Range("E41").Value = "  -2.46%  "
$ws.Range("D42").Value = "0.307"
$ws.Range("E42").Value = "  -3.03%  "
$ws.Range("D43").Value = "4.49"
$ws.Range("E43").Value = "  -3.85%  "
$ws.Range("D44").Value = "1.53"
$ws.Range("E44").Value = "  -4.05%  "
$ws.Range("D45").Value = "1.10"
$ws.Range("E45").Value = "  +0.52%  "
$ws.Range("D46").Value = "2.10"
$ws.Range("E46").Value = "  -4.33%  "
$ws.Range("D47").Value = "135.05"
$ws.Range("E47").Value = "  -4.16%  "
$ws.Range("D48").Value = "3.37"
$ws.Range("E48").Value = "  -2.98%  "
$ws.Range("D49").Value = "0.0718"
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("D50").Value = "0.486"
$ws.Range("E50").Value = "  -6.32%  "
$ws.Range("D51").Value = "0.564"
$ws.Range("E51").Value = "  -1.98%  "
